$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9545444846153259
$ws.Range("B1").Value = 2.200924634933472
$ws.Range("C1").Value = 8.401504516601562
$ws.Range("D1").Value = 1.688021540641785
$ws.Range("E1").Value = 1.394658088684082
